$wb = $excel.ActiveWorkbook

# --- Sheet "Logs": append new row 5 with the second test-mail log entry ---
$ws = $wb.Worksheets.Item("Logs")

$ws.Cells.Item(5, 1).Value = "Kun je 10 dozen schroeven bestellen?"
$ws.Cells.Item(5, 2).Value = "mailmind.test@zohomail.eu"
$ws.Cells.Item(5, 3).Value = "Testmail #2: Kun je 10 dozen schroeven bestellen?"
$ws.Cells.Item(5, 4).Value = "Bestelling / Levering"
$ws.Cells.Item(5, 5).Value = "Geachte klant,`nBedankt voor uw e-mail. Helaas kan ik geen bestellingen plaatsen, maar ik kan u doorverwijzen naar het bestelteam binnen ons bedrijf. Graag ontvang ik de contactgegevens van uw bedrijf, zodat ik de juiste persoon met u in contact kan brengen.`nIk zie uw reactie graag tegemoet.`nMet vriendelijke groet,`n[Naam]`nE-mailassistent"
$ws.Cells.Item(5, 6).Value = "2025-06-29 14:01:39"
$ws.Cells.Item(5, 7).Value = "Ja"
$ws.Cells.Item(5, 8).Value = "Ja"
$ws.Cells.Item(5, 9).Value = "Nee"

# Re-fit the new row's height so no explicit custom height is stored
# (mirrors the untouched rows above it, which rely on the default height).
$ws.Rows.Item(5).AutoFit()

# Extend the four conditional-formatting blocks (D, G, H, I columns) down to row 5
$ws.Range("D2:D4").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("D2:D5"))
$ws.Range("G2:G4").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("G2:G5"))
$ws.Range("H2:H4").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("H2:H5"))
$ws.Range("I2:I4").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("I2:I5"))

# --- Sheet "Dashboard": bump the "Bestelling / Levering" count from 1 to 2 ---
$ws2 = $wb.Worksheets.Item("Dashboard")
$ws2.Range("B3").Value = 2
